# Restructure ontology: remove mfd_hab1=Urban if mfd_areatype=Urban
#
# For every row whose mfd_hab1 (column N) is "Urban" (which, in this
# workbook, is exactly the set of rows where mfd_areatype / column L is
# also "Urban" and habitat_typenumber / column F is "6100"):
#   - habitat_typenumber (F) changes from 6100 to 1000
#   - mfd_hab1 (N) no longer holds "Urban" - the value that used to sit in
#     mfd_hab2 (O, "Roadside") shifts left into mfd_hab1 (N)
#   - mfd_hab2 (O) is removed (there is no mfd_hab3 value to shift into it)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(5, 16, 18, 20, 25, 31, 41, 48, 60, 67)

foreach ($r in $rows) {
    $fCell = $ws.Cells.Item($r, 6)   # F: habitat_typenumber
    $nCell = $ws.Cells.Item($r, 14)  # N: mfd_hab1
    $oCell = $ws.Cells.Item($r, 15)  # O: mfd_hab2

    # Write "1000" as literal text (not a number) into F, matching the
    # existing text-typed habitat_typenumber column: build it via a text
    # formula, then bake the formula's result down to a plain value so the
    # cell ends up a normal literal string cell (no residual formula, no
    # number-format style change).
    $fCell.Formula = '="1000"'
    $fCell.Copy()
    $fCell.PasteSpecial(-4163)

    # mfd_hab2's "Roadside" value moves left into mfd_hab1 ...
    $nCell.Value = "Roadside"
    # ... and mfd_hab2 is dropped entirely.
    $oCell.ClearContents()
}

$excel.CutCopyMode = 0
